# Leave Card update (4/27/2023 4:55 PM):
# Insert a new entry row into Table1 (Sheet1) right above the existing
# 5/1/2023 row (sheet row 107), recording an SL(3-0-00) leave of 3 days
# taken on 4/3,4,20/2023. Every row from 107 downwards shifts down by one,
# and the table grows by one row (A8:K139 -> A8:K140).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift row 107 (and everything below it) down by one row, opening up a
# blank row 107 for the new entry.
$ws.Rows.Item(107).Insert()

# The freshly-inserted row 107 has no formatting yet; pick up the normal
# data-row look (borders/number formats/etc.) from the row right below it
# (the row that used to be 107, now shifted to 108).
$ws.Range("A108:K108").Copy()
$ws.Range("A107:K107").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The table definition still only covers A8:K139; grow it to include the
# new last row (140) so the calculated/structured-reference columns keep
# working for every row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K140"))

# Re-establish the calculated "EARNED " column formula on the new row and
# on the row that now sits at the very bottom of the table.
$ws.Range("G107").Formula = $ws.Range("G108").Formula
$ws.Range("G140").Formula = $ws.Range("G139").Formula

# Fill in the new leave entry: PARTICULARS (SL(3-0-00)), NUMBER OF DAYS (3),
# and REMARKS (4/3,4,20/2023). PERIOD (column A) is left blank, matching
# the edit.
$ws.Range("B107").Value = "SL(3-0-00)"
$ws.Range("D107").Value = 3
$ws.Range("K107").Value = "4/3,4,20/2023"

# Match the author's final on-screen selection.
$ws.Range("H107").Select()
